$d = $word.ActiveDocument

# Locate the UC-50 title paragraph ("Realizar login no aplicativo móvel.")
# without relying on a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*Realizar*login*no aplicativo*") {
        $target = $cand
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Realizar login no aplicativo movel' paragraph"
}

$r = $target.Range
# Exclude the trailing paragraph mark so the paragraph itself (and the
# bookmark that follows the run) is preserved.
$editRange = $d.Range($r.Start, $r.End - 1)

# Replace the three runs ("Realizar " / "login" / " no aplicativo móvel.")
# with a single upper-cased run, and drop the explicit <w:sz w:val="24"/>
# from both the paragraph mark run properties and the run itself (falls
# back to the document default size).
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body><w:p w:rsidR="00C764C2" w:rsidRPr="009946EE" w:rsidRDefault="00C764C2" w:rsidP="00C764C2">' +
'<w:pPr><w:pStyle w:val="Cabealho"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:b/><w:bCs/><w:color w:val="000000" w:themeColor="text1"/></w:rPr>' +
'<w:t>REALIZAR LOGIN NO APLICATIVO MÓVEL.</w:t></w:r>' +
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
'</w:p></w:body>' +
'</w:document></pkg:xmlData></pkg:part></pkg:package>'

$editRange.InsertXML($xml)
